$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '304.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.08%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '31.90'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.40%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.210'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '2.10%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07812'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.77%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.367'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '38.41%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.984'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.99%'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.863'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.68%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9119'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-1.96%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1735'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.36%'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07370'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.78%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08142'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '1.86%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03049'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.40%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09945'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.55%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001526'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1.35%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006193'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-3.16%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.501'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.12%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.242'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.85%'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.74%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1336'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '0.78%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.678'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2.55%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04648'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.19%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.1563'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.30%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001261'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '3.30%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2.67%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '3.70%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0002737'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '47.30%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01796'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '7.16%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04594'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '1.63%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007272'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.40%'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2.74%'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '8.58%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01097'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-6.47%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00006486'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '8.17%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.20%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-57.48%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.009880'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-23.80%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002096'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.20%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001996'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.13%'
